$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Status" dropdown column (F) for the task rows whose status
# progressed: some "In progress"/"Not started" tasks became "Done", and
# some "Not started" tasks became "In progress".
$ws.Range("F3:F14").Value = "Done"
$ws.Range("F26:F31").Value = "In progress"
$ws.Range("F39:F40").Value = "Done"
$ws.Range("F41:F42").Value = "In progress"
$ws.Range("F43:F47").Value = "Done"
$ws.Range("F48:F54").Value = "In progress"

# Update the sheet's current selection to match the saved view.
$ws.Activate()
$ws.Range("F15").Select()
